# Applies the cryptos.xlsx price/volume/ranking update described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Subscript-digit characters used by a couple of very-small-price cells (e.g. 0.0₃0784).
$sub3 = [char]0x2083
$sub6 = [char]0x2086

$ws.Range("D2").Value = '''34.023.22'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.31%  '
$ws.Range("D3").Value = '''1.783.46'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.63%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").Value = '''226.73'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.11%  '
$ws.Range("D6").Value = '''0.550'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.01%  '
$ws.Range("D7").Value = '''0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").Value = '''32.84'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.70%  '
$ws.Range("D9").Value = '''0.286'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.50%  '
$ws.Range("D10").Value = '''0.0713'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.13%  '
$ws.Range("D11").Value = '''0.0934'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.63%  '
$ws.Range("D12").Value = '''2.040.77'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.62%  '
$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").Value = '''11.18'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.75%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '''1.801.59'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.60%  '
$ws.Range("D15").Value = '''33.990.49'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.38%  '
$ws.Range("D16").Value = '''0.620'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.26%  '
$ws.Range("E17").Value = '  -4.47%  '
$ws.Range("D18").Value = '''67.81'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.39%  '
$ws.Range("D19").Value = '''244.58'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.12%  '
$ws.Range("D20").Value = ('''0.0' + $sub3 + '0785')
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.31%  '
$ws.Range("E21").Value = '  +0.28%  '
$ws.Range("D22").Value = '''10.74'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.69%  '
$ws.Range("D23").Value = '''4.09'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.36%  '
$ws.Range("D24").Value = '''2.07'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.68%  '
$ws.Range("D25").Value = '''160.18'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.08%  '
$ws.Range("D26").Value = '''16.31'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.94%  '
$ws.Range("D27").Value = '''7.07'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.23%  '
$ws.Range("E28").Value = '  -2.03%  '
$ws.Range("E29").Value = '  +0.12%  '
$ws.Range("E30").Value = '  +1.86%  '
$ws.Range("D31").Value = '''0.0512'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.55%  '
$ws.Range("D32").Value = '''3.64'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.57%  '
$ws.Range("E33").Value = '  -2.33%  '
$ws.Range("D34").Value = '''1.80'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.14%  '
$ws.Range("D35").Value = '''1.390.73'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.79%  '
$ws.Range("D36").Value = '''0.652'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.50%  '
$ws.Range("E37").Value = '  -2.04%  '
$ws.Range("D38").Value = '''0.0187'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.95%  '
$ws.Range("E39").Value = '  +0.60%  '
$ws.Range("E40").Value = '  +1.21%  '
$ws.Range("D41").Value = '''0.914'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.82%  '
$ws.Range("E42").Value = '  -2.44%  '
$ws.Range("D43").Value = '''78.08'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.94%  '
$ws.Range("B44").Value = 'InjectiveProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D44").Value = '''13.11'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +10.67%  '
$ws.Range("B45").Value = 'BabyDogeCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D45").Value = ('''0.0' + $sub6 + '0138')
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +11.45%  '
$ws.Range("B46").Value = 'WEMIXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D46").Value = '''1.08'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.92%  '
$ws.Range("B47").Value = 'Quant'
$ws.Range("C47").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D47").Value = '''108.26'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.90%  '
$ws.Range("B48").Value = 'Kaspa'
$ws.Range("C48").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D48").Value = '''0.0497'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.22%  '
$ws.Range("B49").Value = 'RocketPoolETH'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D49").Value = '''1.940.37'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.28%  '
$ws.Range("B50").Value = 'FraxShare'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D50").Value = '''5.78'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.47%  '
$ws.Range("D51").Value = '''0.999'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.26%  '
